$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append to the "master-reg_center_machine_devic" sheet.
# Columns: A=regcntr_id, B=machine_id, C=device_id, D=lang_code,
#          E=is_active, F=cr_by, G=cr_dtimes, H=eff_dtimes
$newRows = @(
    @(10002, 10032, 3000176),
    @(10002, 10032, 3000177),
    @(10002, 10032, 3000178),
    @(10002, 10032, 3000179),
    @(10002, 10032, 3000180)
)

$startRow = 157
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Update view: scroll so row 151 is at top, select D157 as active cell.
$ws.Range("D157").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 151

# Switch workbook calculation mode to manual.
$excel.Calculation = -4135
